$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "epurdom"
$ws.Range("B9").Value = "samwise"
$ws.Range("C9").Value = "AMD EPYC 7543 32-Core Processor"
$ws.Range("D9").Value = 60132

$ws.Range("B10").Select()
